$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Combined")
$ws.Cells.Item(2, 2).Value = "BTC-USDT"
$ws.Cells.Item(2, 3).Value = 0.0002329612166132676
$ws.Cells.Item(2, 4).Value = -0.0002329612166132676
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(3, 2).Value = "ETH-USDT"
$ws.Cells.Item(3, 3).Value = 0.0004516079500956745
$ws.Cells.Item(3, 4).Value = -0.0004516079500956745
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(4, 2).Value = "INJ-USDT"
$ws.Cells.Item(4, 3).Value = 0.09479824987846737
$ws.Cells.Item(4, 4).Value = -0.06319883325231734
$ws.Cells.Item(4, 5).Value = 0.03162940074450287
$ws.Cells.Item(5, 2).Value = "AVAX-USDT"
$ws.Cells.Item(5, 3).Value = 0.002406275566672288
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 0.002406333469686622
$ws.Cells.Item(6, 2).Value = "BTC-USDC"
$ws.Cells.Item(6, 3).Value = 0.0002330714387402507
$ws.Cells.Item(6, 4).Value = 0.03099850135063881
$ws.Cells.Item(6, 5).Value = 0.03123164558142476
$ws.Cells.Item(7, 2).Value = "SOL-USDT"
$ws.Cells.Item(7, 3).Value = 0.01711066652627834
$ws.Cells.Item(7, 4).Value = -0.01184584605664992
$ws.Cells.Item(7, 5).Value = 0.005265721469669298
$ws.Cells.Item(8, 2).Value = "FET-USDT"
$ws.Cells.Item(8, 3).Value = 0.02802690582959333
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0.02803476310624866
$ws.Cells.Item(9, 2).Value = "DOGE-USDT"
$ws.Cells.Item(9, 3).Value = 0.01106194690265058
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 0.01106317070471968
$ws.Cells.Item(10, 2).Value = "ETH-USDC"
$ws.Cells.Item(10, 3).Value = 0.0004518957024817334
$ws.Cells.Item(10, 4).Value = 0.01310497537168257
$ws.Cells.Item(10, 5).Value = 0.01355693233735893
$ws.Cells.Item(11, 2).Value = "ADA-USDT"
$ws.Cells.Item(11, 3).Value = 0.01119013911377667
$ws.Cells.Item(11, 4).Value = 0.0003390951246507725
$ws.Cells.Item(11, 5).Value = 0.0115305245201618
$ws.Cells.Item(12, 2).Value = "MATIC-USDT"
$ws.Cells.Item(12, 3).Value = 0.01286504567091071
$ws.Cells.Item(12, 4).Value = 0.02573009134182143
$ws.Cells.Item(12, 5).Value = 0.03860010293360357
$ws.Cells.Item(13, 2).Value = "RUNE-USDT"
$ws.Cells.Item(13, 3).Value = 0.001886685659297907
$ws.Cells.Item(13, 4).Value = 0.01698017093371467
$ws.Cells.Item(13, 5).Value = 0.01886721255800622
$ws.Cells.Item(14, 2).Value = "BONK-USDT"
$ws.Cells.Item(14, 3).Value = 0.06707282991447772
$ws.Cells.Item(14, 4).Value = 0.005589402492885768
$ws.Cells.Item(14, 5).Value = 0.07271100173388557
$ws.Cells.Item(15, 2).Value = "LINK-USDT"
$ws.Cells.Item(15, 3).Value = 0.0006948352892930534
$ws.Cells.Item(15, 4).Value = -0.001389670578586107
$ws.Cells.Item(15, 5).Value = -0.0006948401172873924
$ws.Cells.Item(16, 2).Value = "ATOM-USDT"
$ws.Cells.Item(16, 3).Value = 0.0009400616680599317
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(16, 5).Value = 0.0009400705053024047
$ws.Cells.Item(17, 2).Value = "ETH-BTC"
$ws.Cells.Item(17, 3).Value = 0.001938548027529321
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 0.001938585607942385
$ws.Cells.Item(18, 2).Value = "TIA-USDT"
$ws.Cells.Item(18, 3).Value = 0.09351052785899837
$ws.Cells.Item(18, 4).Value = -0.04517033972849575
$ws.Cells.Item(18, 5).Value = 0.04838543360487339
$ws.Cells.Item(19, 2).Value = "JTO-USDT"
$ws.Cells.Item(19, 3).Value = 0.1384962565276573
$ws.Cells.Item(19, 4).Value = -0.02484785778878504
$ws.Cells.Item(19, 5).Value = 0.113806015810473
$ws.Cells.Item(20, 2).Value = "XRP-USDT"
$ws.Cells.Item(20, 3).Value = 0.001640070194996882
$ws.Cells.Item(20, 4).Value = -0.008200350975020828
$ws.Cells.Item(20, 5).Value = -0.00656038837499836
$ws.Cells.Item(21, 2).Value = "DOT-USDT"
$ws.Cells.Item(21, 3).Value = 0.04787673915881171
$ws.Cells.Item(21, 4).Value = 0.02321296444065043
$ws.Cells.Item(21, 5).Value = 0.07112375533428349
$ws.Cells.Item(22, 2).Value = "BTCUSDT"
$ws.Cells.Item(22, 3).Value = [double]"2.329487350193526e-05"
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = [double]"2.329487892844784e-05"
$ws.Cells.Item(23, 2).Value = "ETHUSDT"
$ws.Cells.Item(23, 3).Value = 0.0004516099896228279
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = 0.0004516120291478658
$ws.Cells.Item(24, 2).Value = "SOLUSDT"
$ws.Cells.Item(24, 3).Value = 0.01316309069370161
$ws.Cells.Item(24, 4).Value = 0
$ws.Cells.Item(24, 5).Value = 0.01316482359137061
$ws.Cells.Item(25, 2).Value = "BTCUSDC"
$ws.Cells.Item(25, 3).Value = [double]"2.330782292589456e-05"
$ws.Cells.Item(25, 4).Value = -0.0001864625833732392
$ws.Cells.Item(25, 5).Value = -0.0001631547984751761
$ws.Cells.Item(26, 2).Value = "INJUSDT"
$ws.Cells.Item(26, 3).Value = 0.0121616033857966
$ws.Cells.Item(26, 4).Value = -0.009729282708626915
$ws.Cells.Item(26, 5).Value = 0.002432616522343035
$ws.Cells.Item(27, 2).Value = "AVAXUSDT"
$ws.Cells.Item(27, 3).Value = 0.024044241404196
$ws.Cells.Item(27, 4).Value = -0.024044241404196
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(28, 2).Value = "ADAUSDT"
$ws.Cells.Item(28, 3).Value = 0.0169520257670773
$ws.Cells.Item(28, 4).Value = -0.0169520257670773
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 2).Value = "TIAUSDT"
$ws.Cells.Item(29, 3).Value = 0.03962624525476333
$ws.Cells.Item(29, 4).Value = -0.00317009962038726
$ws.Cells.Item(29, 5).Value = 0.03647059756281218
$ws.Cells.Item(30, 2).Value = "BONKUSDT"
$ws.Cells.Item(30, 3).Value = 0.05586592178771842
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0.05589714924539742
$ws.Cells.Item(31, 2).Value = "XRPUSDT"
$ws.Cells.Item(31, 3).Value = 0.01640419947506381
$ws.Cells.Item(31, 4).Value = -0.01640419947506381
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 2).Value = "JTOUSDT"
$ws.Cells.Item(32, 3).Value = 0.02038403522361967
$ws.Cells.Item(32, 4).Value = 0.004076807044713072
$ws.Cells.Item(32, 5).Value = 0.02446582939161368
$ws.Cells.Item(33, 2).Value = "RUNEUSDT"
$ws.Cells.Item(33, 3).Value = 0.01886436521411684
$ws.Cells.Item(33, 4).Value = 0
$ws.Cells.Item(33, 5).Value = 0.01886792452830819
$ws.Cells.Item(34, 2).Value = "DOGEUSDT"
$ws.Cells.Item(34, 3).Value = 0.01106194690265058
$ws.Cells.Item(34, 4).Value = -0.02212389380530116
$ws.Cells.Item(34, 5).Value = -0.01106317070471968
$ws.Cells.Item(35, 2).Value = "LINKUSDT"
$ws.Cells.Item(35, 3).Value = 0.00694637399278426
$ws.Cells.Item(35, 4).Value = -0.0208391219783281
$ws.Cells.Item(35, 5).Value = -0.01389371309481689
$ws.Cells.Item(36, 2).Value = "ETHUSDC"
$ws.Cells.Item(36, 3).Value = 0.0004518609894950579
$ws.Cells.Item(36, 4).Value = 0.009489080779190733
$ws.Cells.Item(36, 5).Value = 0.009940986688126605
$ws.Cells.Item(37, 2).Value = "FETUSDT"
$ws.Cells.Item(37, 3).Value = 0.02801120448178964
$ws.Cells.Item(37, 4).Value = -0.01400560224089482
$ws.Cells.Item(37, 5).Value = 0.0140095264780035
$ws.Cells.Item(38, 2).Value = "ETHBTC"
$ws.Cells.Item(38, 3).Value = 0.0193836014731462
$ws.Cells.Item(38, 4).Value = 0
$ws.Cells.Item(38, 5).Value = 0.01938735944163653
$ws.Cells.Item(39, 2).Value = "ATOMUSDT"
$ws.Cells.Item(39, 3).Value = 0.009395847035621744
$ws.Cells.Item(39, 4).Value = -0.0187916940712268
$ws.Cells.Item(39, 5).Value = -0.009396729937976375
$ws.Cells.Item(40, 2).Value = "MATICUSDT"
$ws.Cells.Item(40, 3).Value = 0.01286008230452533
$ws.Cells.Item(40, 4).Value = -0.01286008230452533
$ws.Cells.Item(40, 5).Value = 0
$ws.Cells.Item(41, 2).Value = "DOTUSDT"
$ws.Cells.Item(41, 3).Value = 0.01450747134773605
$ws.Cells.Item(41, 4).Value = -0.01450747134773605
$ws.Cells.Item(41, 5).Value = 0

$ws = $wb.Worksheets.Item("KuCoin")
$ws.Cells.Item(2, 2).Value = "BTC-USDT"
$ws.Cells.Item(2, 3).Value = 0.0002329612166132676
$ws.Cells.Item(2, 4).Value = -0.0002329612166132676
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(3, 2).Value = "ETH-USDT"
$ws.Cells.Item(3, 3).Value = 0.0004516079500956745
$ws.Cells.Item(3, 4).Value = -0.0004516079500956745
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(4, 2).Value = "INJ-USDT"
$ws.Cells.Item(4, 3).Value = 0.09479824987846737
$ws.Cells.Item(4, 4).Value = -0.06319883325231734
$ws.Cells.Item(4, 5).Value = 0.03162940074450287
$ws.Cells.Item(5, 2).Value = "AVAX-USDT"
$ws.Cells.Item(5, 3).Value = 0.002406275566672288
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 0.002406333469686622
$ws.Cells.Item(6, 2).Value = "BTC-USDC"
$ws.Cells.Item(6, 3).Value = 0.0002330714387402507
$ws.Cells.Item(6, 4).Value = 0.03099850135063881
$ws.Cells.Item(6, 5).Value = 0.03123164558142476
$ws.Cells.Item(7, 2).Value = "SOL-USDT"
$ws.Cells.Item(7, 3).Value = 0.01711066652627834
$ws.Cells.Item(7, 4).Value = -0.01184584605664992
$ws.Cells.Item(7, 5).Value = 0.005265721469669298
$ws.Cells.Item(8, 2).Value = "FET-USDT"
$ws.Cells.Item(8, 3).Value = 0.02802690582959333
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0.02803476310624866
$ws.Cells.Item(9, 2).Value = "DOGE-USDT"
$ws.Cells.Item(9, 3).Value = 0.01106194690265058
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 0.01106317070471968
$ws.Cells.Item(10, 2).Value = "ETH-USDC"
$ws.Cells.Item(10, 3).Value = 0.0004518957024817334
$ws.Cells.Item(10, 4).Value = 0.01310497537168257
$ws.Cells.Item(10, 5).Value = 0.01355693233735893
$ws.Cells.Item(11, 2).Value = "ADA-USDT"
$ws.Cells.Item(11, 3).Value = 0.01119013911377667
$ws.Cells.Item(11, 4).Value = 0.0003390951246507725
$ws.Cells.Item(11, 5).Value = 0.0115305245201618
$ws.Cells.Item(12, 2).Value = "MATIC-USDT"
$ws.Cells.Item(12, 3).Value = 0.01286504567091071
$ws.Cells.Item(12, 4).Value = 0.02573009134182143
$ws.Cells.Item(12, 5).Value = 0.03860010293360357
$ws.Cells.Item(13, 2).Value = "RUNE-USDT"
$ws.Cells.Item(13, 3).Value = 0.001886685659297907
$ws.Cells.Item(13, 4).Value = 0.01698017093371467
$ws.Cells.Item(13, 5).Value = 0.01886721255800622
$ws.Cells.Item(14, 2).Value = "BONK-USDT"
$ws.Cells.Item(14, 3).Value = 0.06707282991447772
$ws.Cells.Item(14, 4).Value = 0.005589402492885768
$ws.Cells.Item(14, 5).Value = 0.07271100173388557
$ws.Cells.Item(15, 2).Value = "LINK-USDT"
$ws.Cells.Item(15, 3).Value = 0.0006948352892930534
$ws.Cells.Item(15, 4).Value = -0.001389670578586107
$ws.Cells.Item(15, 5).Value = -0.0006948401172873924
$ws.Cells.Item(16, 2).Value = "ATOM-USDT"
$ws.Cells.Item(16, 3).Value = 0.0009400616680599317
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(16, 5).Value = 0.0009400705053024047
$ws.Cells.Item(17, 2).Value = "ETH-BTC"
$ws.Cells.Item(17, 3).Value = 0.001938548027529321
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 0.001938585607942385
$ws.Cells.Item(18, 2).Value = "TIA-USDT"
$ws.Cells.Item(18, 3).Value = 0.09351052785899837
$ws.Cells.Item(18, 4).Value = -0.04517033972849575
$ws.Cells.Item(18, 5).Value = 0.04838543360487339
$ws.Cells.Item(19, 2).Value = "JTO-USDT"
$ws.Cells.Item(19, 3).Value = 0.1384962565276573
$ws.Cells.Item(19, 4).Value = -0.02484785778878504
$ws.Cells.Item(19, 5).Value = 0.113806015810473
$ws.Cells.Item(20, 2).Value = "XRP-USDT"
$ws.Cells.Item(20, 3).Value = 0.001640070194996882
$ws.Cells.Item(20, 4).Value = -0.008200350975020828
$ws.Cells.Item(20, 5).Value = -0.00656038837499836
$ws.Cells.Item(21, 2).Value = "DOT-USDT"
$ws.Cells.Item(21, 3).Value = 0.04787673915881171
$ws.Cells.Item(21, 4).Value = 0.02321296444065043
$ws.Cells.Item(21, 5).Value = 0.07112375533428349

$ws = $wb.Worksheets.Item("Binance")
$ws.Cells.Item(2, 2).Value = "BTCUSDT"
$ws.Cells.Item(2, 3).Value = [double]"2.329487350193526e-05"
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = [double]"2.329487892844784e-05"
$ws.Cells.Item(3, 2).Value = "ETHUSDT"
$ws.Cells.Item(3, 3).Value = 0.0004516099896228279
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 0.0004516120291478658
$ws.Cells.Item(4, 2).Value = "SOLUSDT"
$ws.Cells.Item(4, 3).Value = 0.01316309069370161
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0.01316482359137061
$ws.Cells.Item(5, 2).Value = "BTCUSDC"
$ws.Cells.Item(5, 3).Value = [double]"2.330782292589456e-05"
$ws.Cells.Item(5, 4).Value = -0.0001864625833732392
$ws.Cells.Item(5, 5).Value = -0.0001631547984751761
$ws.Cells.Item(6, 2).Value = "INJUSDT"
$ws.Cells.Item(6, 3).Value = 0.0121616033857966
$ws.Cells.Item(6, 4).Value = -0.009729282708626915
$ws.Cells.Item(6, 5).Value = 0.002432616522343035
$ws.Cells.Item(7, 2).Value = "AVAXUSDT"
$ws.Cells.Item(7, 3).Value = 0.024044241404196
$ws.Cells.Item(7, 4).Value = -0.024044241404196
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(8, 2).Value = "ADAUSDT"
$ws.Cells.Item(8, 3).Value = 0.0169520257670773
$ws.Cells.Item(8, 4).Value = -0.0169520257670773
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(9, 2).Value = "TIAUSDT"
$ws.Cells.Item(9, 3).Value = 0.03962624525476333
$ws.Cells.Item(9, 4).Value = -0.00317009962038726
$ws.Cells.Item(9, 5).Value = 0.03647059756281218
$ws.Cells.Item(10, 2).Value = "BONKUSDT"
$ws.Cells.Item(10, 3).Value = 0.05586592178771842
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 0.05589714924539742
$ws.Cells.Item(11, 2).Value = "XRPUSDT"
$ws.Cells.Item(11, 3).Value = 0.01640419947506381
$ws.Cells.Item(11, 4).Value = -0.01640419947506381
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(12, 2).Value = "JTOUSDT"
$ws.Cells.Item(12, 3).Value = 0.02038403522361967
$ws.Cells.Item(12, 4).Value = 0.004076807044713072
$ws.Cells.Item(12, 5).Value = 0.02446582939161368
$ws.Cells.Item(13, 2).Value = "RUNEUSDT"
$ws.Cells.Item(13, 3).Value = 0.01886436521411684
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0.01886792452830819
$ws.Cells.Item(14, 2).Value = "DOGEUSDT"
$ws.Cells.Item(14, 3).Value = 0.01106194690265058
$ws.Cells.Item(14, 4).Value = -0.02212389380530116
$ws.Cells.Item(14, 5).Value = -0.01106317070471968
$ws.Cells.Item(15, 2).Value = "LINKUSDT"
$ws.Cells.Item(15, 3).Value = 0.00694637399278426
$ws.Cells.Item(15, 4).Value = -0.0208391219783281
$ws.Cells.Item(15, 5).Value = -0.01389371309481689
$ws.Cells.Item(16, 2).Value = "ETHUSDC"
$ws.Cells.Item(16, 3).Value = 0.0004518609894950579
$ws.Cells.Item(16, 4).Value = 0.009489080779190733
$ws.Cells.Item(16, 5).Value = 0.009940986688126605
$ws.Cells.Item(17, 2).Value = "FETUSDT"
$ws.Cells.Item(17, 3).Value = 0.02801120448178964
$ws.Cells.Item(17, 4).Value = -0.01400560224089482
$ws.Cells.Item(17, 5).Value = 0.0140095264780035
$ws.Cells.Item(18, 2).Value = "ETHBTC"
$ws.Cells.Item(18, 3).Value = 0.0193836014731462
$ws.Cells.Item(18, 4).Value = 0
$ws.Cells.Item(18, 5).Value = 0.01938735944163653
$ws.Cells.Item(19, 2).Value = "ATOMUSDT"
$ws.Cells.Item(19, 3).Value = 0.009395847035621744
$ws.Cells.Item(19, 4).Value = -0.0187916940712268
$ws.Cells.Item(19, 5).Value = -0.009396729937976375
$ws.Cells.Item(20, 2).Value = "MATICUSDT"
$ws.Cells.Item(20, 3).Value = 0.01286008230452533
$ws.Cells.Item(20, 4).Value = -0.01286008230452533
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(21, 2).Value = "DOTUSDT"
$ws.Cells.Item(21, 3).Value = 0.01450747134773605
$ws.Cells.Item(21, 4).Value = -0.01450747134773605
$ws.Cells.Item(21, 5).Value = 0
